# edit.ps1 -- apply the "added modal select cities 80%" change to notes.docx
#
# Summary of changes (per the unified OOXML diff):
#   1. The empty paragraph right after the "...RouteMaster.git" hyperlink
#      paragraph (originally just an empty <w:p>) is replaced by two new
#      paragraphs:
#         - a bold "Domain:" heading paragraph (preceded by a line break run)
#         - a paragraph containing "https://dash.cloudflare.com/" followed
#           by a line break run
#   2. The stale <w:lastRenderedPageBreak/> marker inside the separator-dash
#      paragraph right before "КРИВОЙ РОГ" (the 4th "----...----" separator
#      paragraph in the document) is removed; nothing else about that
#      paragraph changes.
#   3. The section's top page margin shrinks from 1134 twips (56.7pt) to
#      568 twips (28.4pt).

$d = $word.ActiveDocument
$CR = [string][char]13

# ---------------------------------------------------------------------
# 1) Replace the empty paragraph that immediately follows the
#    ".../RouteMaster.git" hyperlink paragraph with the new
#    "Domain:" / cloudflare block.
# ---------------------------------------------------------------------
$targetParaIndex = 0
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq $CR) {
        $prevText = $d.Paragraphs($i - 1).Range.Text
        if ($prevText -like "*RouteMaster.git*") {
            $targetParaIndex = $i
            break
        }
    }
}

if ($targetParaIndex -eq 0) {
    throw "Could not locate the empty placeholder paragraph to replace"
}

$target = $d.Paragraphs($targetParaIndex).Range

$domainXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:pPr><w:rPr><w:rStyle w:val="a3"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="none"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="a3"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="none"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rStyle w:val="a3"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="auto"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="none"/><w:lang w:val="en-US"/></w:rPr><w:t>Domain:</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:rPr><w:rStyle w:val="a3"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="a3"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>https://dash.cloudflare.com/</w:t></w:r><w:r><w:rPr><w:rStyle w:val="a3"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($domainXml)

# ---------------------------------------------------------------------
# 2) Strip the stale <w:lastRenderedPageBreak/> from the 4th separator
#    ("----...----") paragraph -- leave everything else about that
#    paragraph (including its identity / rsid attributes) untouched by
#    re-writing only the run content, not the paragraph mark itself.
# ---------------------------------------------------------------------
$dashes = "".PadLeft(74, "-")

$dashIndices = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq $dashes) {
        $dashIndices += $i
    }
}

if ($dashIndices.Count -lt 4) {
    throw "Expected at least 4 separator paragraphs, found $($dashIndices.Count)"
}

$sepIdx = $dashIndices[3]
$sepPara = $d.Paragraphs($sepIdx)
$full = $sepPara.Range
$sub = $d.Range($full.Start, $full.End - 1)

$dashXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r w:rsidRPr="00221FE0"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>' + $dashes + '</w:t></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$sub.InsertXML($dashXml)

# ---------------------------------------------------------------------
# 3) Shrink the top page margin from 1134 twips (56.7pt) to 568 twips
#    (28.4pt).
# ---------------------------------------------------------------------
$d.PageSetup.TopMargin = 28.4

Write-Output "edit complete"
